$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "16456234"
$ws.Range("B3").NumberFormat = "General"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "16456242"
$ws.Range("B4").NumberFormat = "General"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "16456243"
$ws.Range("B5").NumberFormat = "General"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "16456245"
$ws.Range("B6").NumberFormat = "General"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "16456246"
$ws.Range("B7").NumberFormat = "General"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "16456244"
$ws.Range("B8").NumberFormat = "General"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "16456254"
$ws.Range("B9").NumberFormat = "General"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "16456253"
$ws.Range("B10").NumberFormat = "General"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "16456250"
$ws.Range("B11").NumberFormat = "General"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "16456256"
$ws.Range("B12").NumberFormat = "General"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "16456236"
$ws.Range("B13").NumberFormat = "General"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "16456233"
$ws.Range("B14").NumberFormat = "General"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "16457945"
$ws.Range("B15").NumberFormat = "General"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "16456218"
$ws.Range("B16").NumberFormat = "General"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "16456221"
$ws.Range("B17").NumberFormat = "General"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "16456227"
$ws.Range("B18").NumberFormat = "General"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "16456205"
$ws.Range("B19").NumberFormat = "General"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "16456208"
$ws.Range("B20").NumberFormat = "General"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "16456212"
$ws.Range("B21").NumberFormat = "General"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "16456261"
$ws.Range("B22").NumberFormat = "General"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "16456712"
$ws.Range("B23").NumberFormat = "General"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "16456709"
$ws.Range("B24").NumberFormat = "General"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "16456750"
$ws.Range("B25").NumberFormat = "General"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "16456749"
$ws.Range("B26").NumberFormat = "General"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "16456754"
$ws.Range("B27").NumberFormat = "General"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "16457272"
$ws.Range("B28").NumberFormat = "General"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "16456748"
$ws.Range("B29").NumberFormat = "General"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "16456735"
$ws.Range("B30").NumberFormat = "General"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "16456493"
$ws.Range("B31").NumberFormat = "General"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "16456758"
$ws.Range("B32").NumberFormat = "General"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "16456756"
$ws.Range("B33").NumberFormat = "General"
$ws.Range("B34").Value = "LP024176"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "16456734"
$ws.Range("B35").NumberFormat = "General"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "16456733"
$ws.Range("B36").NumberFormat = "General"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "16456731"
$ws.Range("B37").NumberFormat = "General"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "16456764"
$ws.Range("B38").NumberFormat = "General"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "16456765"
$ws.Range("B39").NumberFormat = "General"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "16456701"
$ws.Range("B40").NumberFormat = "General"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "16456706"
$ws.Range("B41").NumberFormat = "General"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "16456708"
$ws.Range("B42").NumberFormat = "General"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "16457230"
$ws.Range("B43").NumberFormat = "General"
$ws.Range("B44").Value = "NI729519"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "16457262"
$ws.Range("B45").NumberFormat = "General"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "16457268"
$ws.Range("B46").NumberFormat = "General"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "16457273"
$ws.Range("B47").NumberFormat = "General"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "16457265"
$ws.Range("B48").NumberFormat = "General"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "16457270"
$ws.Range("B49").NumberFormat = "General"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "16457267"
$ws.Range("B50").NumberFormat = "General"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "16457279"
$ws.Range("B51").NumberFormat = "General"
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = "16457237"
$ws.Range("B53").NumberFormat = "General"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = "16459074"
$ws.Range("B54").NumberFormat = "General"
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = "16457259"
$ws.Range("B55").NumberFormat = "General"
$ws.Range("B56").Value = "NI729522"
$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = "16457261"
$ws.Range("B57").NumberFormat = "General"
$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = "16457264"
$ws.Range("B58").NumberFormat = "General"
$ws.Range("B59").NumberFormat = "@"
$ws.Range("B59").Value = "16457962"
$ws.Range("B59").NumberFormat = "General"
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = "16457957"
$ws.Range("B60").NumberFormat = "General"
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = "16457970"
$ws.Range("B61").NumberFormat = "General"
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = "16457965"
$ws.Range("B62").NumberFormat = "General"
$ws.Range("B63").NumberFormat = "@"
$ws.Range("B63").Value = "16457967"
$ws.Range("B63").NumberFormat = "General"
$ws.Range("B64").NumberFormat = "@"
$ws.Range("B64").Value = "16456203"
$ws.Range("B64").NumberFormat = "General"
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = "16457943"
$ws.Range("B65").NumberFormat = "General"
$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = "16457974"
$ws.Range("B66").NumberFormat = "General"
$ws.Range("B67").NumberFormat = "@"
$ws.Range("B67").Value = "16457975"
$ws.Range("B67").NumberFormat = "General"
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = "16457973"
$ws.Range("B68").NumberFormat = "General"
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = "16457954"
$ws.Range("B69").NumberFormat = "General"
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = "16457953"
$ws.Range("B70").NumberFormat = "General"
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = "16459058"
$ws.Range("B71").NumberFormat = "General"
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = "16459059"
$ws.Range("B72").NumberFormat = "General"
$ws.Range("B73").NumberFormat = "@"
$ws.Range("B73").Value = "16459057"
$ws.Range("B73").NumberFormat = "General"
$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = "16459061"
$ws.Range("B74").NumberFormat = "General"
$ws.Range("B75").NumberFormat = "@"
$ws.Range("B75").Value = "16459062"
$ws.Range("B75").NumberFormat = "General"
$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = "16457238"
$ws.Range("B76").NumberFormat = "General"
$ws.Range("B77").NumberFormat = "@"
$ws.Range("B77").Value = "16459078"
$ws.Range("B77").NumberFormat = "General"
$ws.Range("B78").NumberFormat = "@"
$ws.Range("B78").Value = "16459064"
$ws.Range("B78").NumberFormat = "General"
$ws.Range("B79").NumberFormat = "@"
$ws.Range("B79").Value = "16459077"
$ws.Range("B79").NumberFormat = "General"
$ws.Range("B80").NumberFormat = "@"
$ws.Range("B80").Value = "16459084"
$ws.Range("B80").NumberFormat = "General"
$ws.Range("B81").NumberFormat = "@"
$ws.Range("B81").Value = "16459082"
$ws.Range("B81").NumberFormat = "General"
$ws.Range("B82").NumberFormat = "@"
$ws.Range("B82").Value = "16459072"
$ws.Range("B82").NumberFormat = "General"
$ws.Range("B83").NumberFormat = "@"
$ws.Range("B83").Value = "16459526"
$ws.Range("B83").NumberFormat = "General"
$ws.Range("B84").NumberFormat = "@"
$ws.Range("B84").Value = "16459522"
$ws.Range("B84").NumberFormat = "General"
$ws.Range("B85").NumberFormat = "@"
$ws.Range("B85").Value = "16459523"
$ws.Range("B85").NumberFormat = "General"
$ws.Range("B86").NumberFormat = "@"
$ws.Range("B86").Value = "16459521"
$ws.Range("B86").NumberFormat = "General"
$ws.Range("B87").NumberFormat = "@"
$ws.Range("B87").Value = "16459520"
$ws.Range("B87").NumberFormat = "General"
$ws.Range("B88").NumberFormat = "@"
$ws.Range("B88").Value = "16459512"
$ws.Range("B88").NumberFormat = "General"
$ws.Range("B89").NumberFormat = "@"
$ws.Range("B89").Value = "16459515"
$ws.Range("B89").NumberFormat = "General"
$ws.Range("B90").NumberFormat = "@"
$ws.Range("B90").Value = "16459518"
$ws.Range("B90").NumberFormat = "General"
$ws.Range("B91").NumberFormat = "@"
$ws.Range("B91").Value = "16459516"
$ws.Range("B91").NumberFormat = "General"
$ws.Range("B92").NumberFormat = "@"
$ws.Range("B92").Value = "16459514"
$ws.Range("B92").NumberFormat = "General"
$ws.Range("B93").NumberFormat = "@"
$ws.Range("B93").Value = "16459549"
$ws.Range("B93").NumberFormat = "General"
$ws.Range("B94").NumberFormat = "@"
$ws.Range("B94").Value = "16459547"
$ws.Range("B94").NumberFormat = "General"
$ws.Range("B95").NumberFormat = "@"
$ws.Range("B95").Value = "16459548"
$ws.Range("B95").NumberFormat = "General"
$ws.Range("B96").NumberFormat = "@"
$ws.Range("B96").Value = "16459546"
$ws.Range("B96").NumberFormat = "General"
$ws.Range("B97").NumberFormat = "@"
$ws.Range("B97").Value = "16459543"
$ws.Range("B97").NumberFormat = "General"
$ws.Range("B98").NumberFormat = "@"
$ws.Range("B98").Value = "16459542"
$ws.Range("B98").NumberFormat = "General"
$ws.Range("B99").NumberFormat = "@"
$ws.Range("B99").Value = "16459541"
$ws.Range("B99").NumberFormat = "General"
$ws.Range("B100").NumberFormat = "@"
$ws.Range("B100").Value = "16459069"
$ws.Range("B100").NumberFormat = "General"
